# "Generate Report for Archive"
# The localization-status report was regenerated: the "Status" value for the
# two language jobs flips from "Ready for handoff" to "In Translation", and
# the "Status" column on each sheet narrows to fit the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status text wherever it appears ---------------------------
# Overview sheet keeps one status column per locale (E = zh-cn, F = de-de).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Each locale sheet has its own "Status" column (column C).
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Re-fit the Status column widths to the new, shorter text -------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
